$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new attendance record row (row 4), mirroring the existing rows' pattern
$ws.Range("A4").Value = "wambua"
$ws.Range("B4").Value = "2024-10-23 12:20:33"
